$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.348624516477416
$ws.Range("E2").Value = 6.253707197847636

$ws.Range("C3").Value = -7.921319741078669
$ws.Range("E3").Value = -14.35806537048445

$ws.Range("C4").Value = 0.2827737265587604
$ws.Range("E4").Value = 0.646969966709543

$ws.Range("C5").Value = 4.453626355171969
$ws.Range("E5").Value = 6.136355062499943

$ws.Range("C6").Value = 1.477633171193093
$ws.Range("E6").Value = 2.015050062499957

$ws.Range("C7").Value = -0.2005250704869344
$ws.Range("E7").Value = 0.2533625145866525

$ws.Range("C8").Value = 2.234527904461125
$ws.Range("E8").Value = 3.31205545981732

$ws.Range("C9").Value = 1.31172787261824
$ws.Range("E9").Value = 1.093673275363716

$ws.Range("C10").Value = 1.784808447869168
$ws.Range("E10").Value = 2.687934870329323

$ws.Range("C11").Value = 1.874682902292824
$ws.Range("E11").Value = 2.445693358388845

$ws.Range("C12").Value = 2.159589514946769
$ws.Range("E12").Value = 1.194430031759008

$ws.Range("C13").Value = 0.801449343663907
$ws.Range("E13").Value = 1.609625625600009

$ws.Range("C14").Value = -3.107661574595755
$ws.Range("E14").Value = -8.513835774399992

$ws.Range("C15").Value = 0.6949587062036411
$ws.Range("E15").Value = 4.530463903052695

$ws.Range("C16").Value = 2.016988966764255
$ws.Range("E16").Value = 1.492291178243965

$ws.Range("C17").Value = -0.1519071329076249
$ws.Range("E17").Value = 0.3688835244738842

$ws.Range("C18").Value = -0.01173401322185352
$ws.Range("E18").Value = 1.025808301409614

$ws.Range("C19").Value = 0.1328390304517146
$ws.Range("E19").Value = 0.2797831653477356
